# "began building out SOAP examples"
#
# Geocode2 already sketches out a second (json) example row; bring the
# same data over to the Geocode sheet (row 2's numbers become quoted
# text, and a matching row 3 is added), then make Geocode the active
# sheet again while Geocode2's whole range stays selected.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Geocode
$ws2 = $wb.Worksheets.Item(2)   # Geocode2

# --- Geocode (sheet1): re-key B2/F2 as quoted text, add row 3 --------------
$ws2.Range("B2").Copy($ws1.Range("B2"))
$ws2.Range("F2").Copy($ws1.Range("F2"))
$ws2.Range("A3:K3").Copy($ws1.Range("A3:K3"))

# --- View/selection bookkeeping --------------------------------------------
# Geocode2 is no longer the active tab; select its whole used range.
$ws2.Range("A1:XFD3").Select()

# Geocode becomes the active tab, with B4 selected.
$ws1.Select()
$ws1.Range("B4").Select()
